$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.653.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.501.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.39"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.638"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.33%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000272"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.26"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.060.90"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.497.89"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.10"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "65.674.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.998"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "86.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.26"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.81"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -8.59%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.04"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.35"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.49"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "607.77"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -11.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.67"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "59.54"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +10.29%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0792"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0792"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.371.55"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.380"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.68%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.42%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.25"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.72"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.46"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.06%  "
